# got all data for november
# - Row 113 (ticker "UPS") is removed from Sheet1 entirely; all rows below it
#   shift up by one (123 rows -> 122 rows). This also removes the now-unused
#   "UPS" shared string.
# - A few existing counts on Sheet1 were corrected: I (row 53) 18->17,
#   the paired puts row (row 54) 28->27, and the OTM puts row (row 78) 4->3.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Delete the whole row for ticker "UPS" (row 113), shifting everything below up.
$ws1.Rows.Item(113).Delete()

# Correct a handful of counts that were updated with new November data.
$ws1.Cells.Item(53, 3).Value = 17
$ws1.Cells.Item(54, 3).Value = 27
$ws1.Cells.Item(78, 3).Value = 3
